$d = $word.ActiveDocument
$tab = [char]9

# --- Change 1: remove the trailing "  <tabs>  December 1, 2022" after
#     "...Linear Models" on the title line (sheet 1 has no date). The run
#     holding "odels" is left untouched; the space/tab/date runs after it
#     are deleted outright so no text content remains in their place.
$rng1 = $d.Content
$search1 = " $tab$tab$tab${tab}December 1, 2022"
$found1 = $rng1.Find.Execute($search1)
if ($found1) {
    $rng1.Delete()
}

# --- Change 2: merge "), e." into the following run so the text reads
#     "...test_results), e.g., Exponential Weighted..." instead of having
#     "), e." in its own run just before <w:proofErr w:type="spellEnd"/>.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("), e.")
if ($found2) {
    $rng2.Delete()
}

$rng3 = $d.Content
$oldText = "g., Exponential Weighted Moving Averages (EWMA, included in python package "
$newText = "), e.g., Exponential Weighted Moving Averages (EWMA, included in python package "
$found3 = $rng3.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
